# Update "想去人数" (F column) figures across all sheets to match the
# refreshed output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 855
$ws.Range("F3").Value = 13755
$ws.Range("F4").Value = 13546
$ws.Range("F6").Value = 806
$ws.Range("F7").Value = 41
$ws.Range("F8").Value = 593
$ws.Range("F9").Value = 81
$ws.Range("F11").Value = 52
$ws.Range("F12").Value = 757
$ws.Range("F14").Value = 92
$ws.Range("F15").Value = 88
$ws.Range("F17").Value = 116
$ws.Range("F19").Value = 519
$ws.Range("F20").Value = 428
$ws.Range("F21").Value = 389
$ws.Range("F23").Value = 259
$ws.Range("F24").Value = 827
$ws.Range("F25").Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 1482
$ws.Range("F11").Value = 64

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 221
$ws.Range("F3").Value = 104

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 221
$ws.Range("F3").Value = 855
$ws.Range("F4").Value = 13755
$ws.Range("F5").Value = 13546
$ws.Range("F7").Value = 806
$ws.Range("F8").Value = 41
$ws.Range("F9").Value = 593
$ws.Range("F10").Value = 81
$ws.Range("F12").Value = 52
$ws.Range("F13").Value = 757
$ws.Range("F17").Value = 92
$ws.Range("F18").Value = 88
$ws.Range("F20").Value = 116
$ws.Range("F24").Value = 104
$ws.Range("F25").Value = 104
$ws.Range("F26").Value = 519
$ws.Range("F27").Value = 428
$ws.Range("F28").Value = 389
$ws.Range("F30").Value = 259
$ws.Range("F31").Value = 827
$ws.Range("F33").Value = 1482
$ws.Range("F37").Value = 80
$ws.Range("F38").Value = 64
